# "merged from eroots devel"
#
# Adds a new "Qgen" worksheet (generator reactive-power results, in MVAr)
# after the existing "Qbranch" sheet, fills in its data, and makes it the
# active/selected sheet - matching the upstream commit that introduced the
# Qgen results tab alongside Vabs / Vang / Pbranch / Qbranch.

$wb = $excel.ActiveWorkbook

# --- Re-settle the selection on the pre-existing sheets -------------------
# (the original file had a stray multi-range selection left over on Vabs /
# Vang / Qbranch, and Pbranch was the tabSelected sheet; normalise all of
# that back to a simple single-cell selection before handing focus to the
# new sheet.)
$wsVabs = $wb.Worksheets.Item("Vabs")
$wsVabs.Activate()
$wsVabs.Range("A1").Select() | Out-Null

$wsVang = $wb.Worksheets.Item("Vang")
$wsVang.Activate()
$wsVang.Range("A1").Select() | Out-Null

$wsPbranch = $wb.Worksheets.Item("Pbranch")
$wsPbranch.Activate()
$wsPbranch.Range("A2").Select() | Out-Null

$wsQbranch = $wb.Worksheets.Item("Qbranch")
$wsQbranch.Activate()
$wsQbranch.Range("A1").Select() | Out-Null

# --- Add the new "Qgen" sheet after "Qbranch" ------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQgen = $wb.Worksheets.Add($null, $lastSheet)
$wsQgen.Name = "Qgen"

# Header
$wsQgen.Range("B1").Value = "MVAr"

# Bus # / generator reactive power (MVAr) results
$wsQgen.Range("A2").Value = 1
$wsQgen.Range("B2").Value = -20.4176

$wsQgen.Range("A3").Value = 2
$wsQgen.Range("B3").Value = 56.0691

$wsQgen.Range("A4").Value = 5
$wsQgen.Range("B4").Value = 35.6585

$wsQgen.Range("A5").Value = 8
$wsQgen.Range("B5").Value = 36.111

$wsQgen.Range("A6").Value = 11
$wsQgen.Range("B6").Value = 16.0574

$wsQgen.Range("A7").Value = 13
$wsQgen.Range("B7").Value = 10.4507

# Make Qgen the active sheet / selection, as it is the newly added tab
$wsQgen.Activate()
$excel.ActiveWindow.Zoom = 90
$wsQgen.Range("A8").Select() | Out-Null
